# "bug fix in template": the "Growth by Tax" sheet had stale/incorrect
# hard-coded growth-rate values in column E ("Tax Base" / "Five Year Plan").
# Clear them out (keeping the cell formatting) across every data block on
# that sheet. Everything downstream (the "Growth by Year" sheet, which
# pulls these via formulas) recalculates on its own.

$wb = $excel.ActiveWorkbook

$growthByTax = $wb.Worksheets.Item("Growth by Tax")
$growthByYear = $wb.Worksheets.Item("Growth by Year")
$revenueByFY = $wb.Worksheets.Item("Revenue by FY")

# Clear the stale column-E values (5-row blocks, separated by blank/header rows)
$growthByTax.Range("E9:E13").ClearContents()
$growthByTax.Range("E16:E20").ClearContents()
$growthByTax.Range("E23:E27").ClearContents()
$growthByTax.Range("E30:E34").ClearContents()
$growthByTax.Range("E37:E41").ClearContents()
$growthByTax.Range("E44:E48").ClearContents()
$growthByTax.Range("E51:E55").ClearContents()

# Leave "Revenue by FY" selected where it was, just no longer the active tab
$revenueByFY.Activate()
$revenueByFY.Range("A65").Select()

# Make "Growth by Year" selection land on E9 before switching away
$growthByYear.Activate()
$growthByYear.Range("E9").Select()

# "Growth by Tax" becomes the active/selected tab, with E9 selected
$growthByTax.Activate()
$growthByTax.Range("E9").Select()
